# Atualizacao de bases das ligas, do dia: 01-04-2024 as 22:23
# - Row 186/187: the two matches swap places (their stats travel with them).
# - Five upcoming fixtures are appended as rows 217:221.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 186 <-> row 187: swap which match/stat-line sits on which row ---
$ws.Cells.Item(186,6).Value2  = "Unam Pumas U23"   # F186 HomeTeam
$ws.Cells.Item(186,7).Value2  = "Tijuana U23"      # G186 AwayTeam
$ws.Cells.Item(187,6).Value2  = "Monterrey U23"    # F187 HomeTeam
$ws.Cells.Item(187,7).Value2  = "Mazatlan FC U23"  # G187 AwayTeam

$ws.Cells.Item(186,2).Value2 = 7648957
$ws.Cells.Item(187,2).Value2 = 7648958
$ws.Cells.Item(186,8).Value2 = 2
$ws.Cells.Item(187,8).Value2 = 4
$ws.Cells.Item(186,9).Value2 = 0
$ws.Cells.Item(187,9).Value2 = 3
$ws.Cells.Item(186,11).Value2 = 1.666
$ws.Cells.Item(187,11).Value2 = 2.375
$ws.Cells.Item(186,12).Value2 = 3.5
$ws.Cells.Item(187,12).Value2 = 3.1
$ws.Cells.Item(186,13).Value2 = 4.2
$ws.Cells.Item(187,13).Value2 = 2.75
$ws.Cells.Item(186,14).Value2 = 1.533
$ws.Cells.Item(187,14).Value2 = 2.375
$ws.Cells.Item(186,15).Value2 = 4.333
$ws.Cells.Item(187,15).Value2 = 3.4
$ws.Cells.Item(186,16).Value2 = 6
$ws.Cells.Item(187,16).Value2 = 3
$ws.Cells.Item(186,17).Value2 = -1.25
$ws.Cells.Item(187,17).Value2 = -0.25
$ws.Cells.Item(186,18).Value2 = 2.025
$ws.Cells.Item(187,18).Value2 = 2
$ws.Cells.Item(186,19).Value2 = 1.775
$ws.Cells.Item(187,19).Value2 = 1.8
$ws.Cells.Item(186,20).Value2 = 2.75
$ws.Cells.Item(187,20).Value2 = 2.75
$ws.Cells.Item(186,21).Value2 = 1.775
$ws.Cells.Item(187,21).Value2 = 1.95
$ws.Cells.Item(186,22).Value2 = 2.025
$ws.Cells.Item(187,22).Value2 = 1.85
$ws.Cells.Item(186,23).Value2 = 0.5329999999999999
$ws.Cells.Item(187,23).Value2 = 1.375
$ws.Cells.Item(186,24).Value2 = -1
$ws.Cells.Item(187,24).Value2 = -1
$ws.Cells.Item(186,25).Value2 = -1
$ws.Cells.Item(187,25).Value2 = -1
$ws.Cells.Item(186,26).Value2 = 1.025
$ws.Cells.Item(187,26).Value2 = 1
$ws.Cells.Item(186,27).Value2 = -1
$ws.Cells.Item(187,27).Value2 = -1
$ws.Cells.Item(186,28).Value2 = -1
$ws.Cells.Item(187,28).Value2 = 0.95
$ws.Cells.Item(186,29).Value2 = 1.025
$ws.Cells.Item(187,29).Value2 = -1

# --- Append 5 new upcoming-fixture rows (217:221) ---
# row 217
$ws.Cells.Item(217,1).Value2 = 215
$ws.Cells.Item(217,2).Value2 = 7640652
$ws.Cells.Item(217,3).Value2 = "Mexico Liga de Expansion"
$ws.Cells.Item(217,4).Value2 = "Mexico Liga de Expansion"
$ws.Cells.Item(217,5).Value2 = 45384.83680555555
$ws.Cells.Item(217,6).Value2 = "Tapatio"
$ws.Cells.Item(217,7).Value2 = "Atlante"
$ws.Cells.Item(217,11).Value2 = 2.6
$ws.Cells.Item(217,12).Value2 = 3.3
$ws.Cells.Item(217,13).Value2 = 2.5
$ws.Cells.Item(217,14).Value2 = 2.7
$ws.Cells.Item(217,15).Value2 = 3.3
$ws.Cells.Item(217,16).Value2 = 2.4
$ws.Cells.Item(217,17).Value2 = 0
$ws.Cells.Item(217,18).Value2 = 2.025
$ws.Cells.Item(217,19).Value2 = 1.775
$ws.Cells.Item(217,20).Value2 = 2.5
$ws.Cells.Item(217,21).Value2 = 1.975
$ws.Cells.Item(217,22).Value2 = 1.825
$ws.Cells.Item(217,23).Value2 = 0
$ws.Cells.Item(217,24).Value2 = 0
$ws.Cells.Item(217,25).Value2 = 0
$ws.Cells.Item(217,26).Value2 = 0
$ws.Cells.Item(217,27).Value2 = 0
# row 218
$ws.Cells.Item(218,1).Value2 = 216
$ws.Cells.Item(218,2).Value2 = 7641717
$ws.Cells.Item(218,3).Value2 = "Mexico Liga de Expansion"
$ws.Cells.Item(218,4).Value2 = "Mexico Liga de Expansion"
$ws.Cells.Item(218,5).Value2 = 45384.92013888889
$ws.Cells.Item(218,6).Value2 = "Cimarrones de Sonora FC"
$ws.Cells.Item(218,7).Value2 = "Club Atletico La Paz"
$ws.Cells.Item(218,11).Value2 = 2.15
$ws.Cells.Item(218,12).Value2 = 3.25
$ws.Cells.Item(218,13).Value2 = 3.25
$ws.Cells.Item(218,14).Value2 = 2
$ws.Cells.Item(218,15).Value2 = 3.4
$ws.Cells.Item(218,16).Value2 = 3.6
$ws.Cells.Item(218,17).Value2 = -0.5
$ws.Cells.Item(218,18).Value2 = 2
$ws.Cells.Item(218,19).Value2 = 1.8
$ws.Cells.Item(218,20).Value2 = 2.5
$ws.Cells.Item(218,21).Value2 = 2
$ws.Cells.Item(218,22).Value2 = 1.8
$ws.Cells.Item(218,23).Value2 = 0
$ws.Cells.Item(218,24).Value2 = 0
$ws.Cells.Item(218,25).Value2 = 0
$ws.Cells.Item(218,26).Value2 = 0
$ws.Cells.Item(218,27).Value2 = 0
# row 219
$ws.Cells.Item(219,1).Value2 = 217
$ws.Cells.Item(219,2).Value2 = 7641718
$ws.Cells.Item(219,3).Value2 = "Mexico Liga de Expansion"
$ws.Cells.Item(219,4).Value2 = "Mexico Liga de Expansion"
$ws.Cells.Item(219,5).Value2 = 45385.00347222222
$ws.Cells.Item(219,6).Value2 = "Dorados"
$ws.Cells.Item(219,7).Value2 = "Oaxaca"
$ws.Cells.Item(219,11).Value2 = 2.25
$ws.Cells.Item(219,12).Value2 = 3.25
$ws.Cells.Item(219,13).Value2 = 3
$ws.Cells.Item(219,14).Value2 = 2.2
$ws.Cells.Item(219,15).Value2 = 3.25
$ws.Cells.Item(219,16).Value2 = 3.1
$ws.Cells.Item(219,17).Value2 = -0.25
$ws.Cells.Item(219,18).Value2 = 1.925
$ws.Cells.Item(219,19).Value2 = 1.875
$ws.Cells.Item(219,20).Value2 = 2.5
$ws.Cells.Item(219,21).Value2 = 2
$ws.Cells.Item(219,22).Value2 = 1.8
$ws.Cells.Item(219,23).Value2 = 0
$ws.Cells.Item(219,24).Value2 = 0
$ws.Cells.Item(219,25).Value2 = 0
$ws.Cells.Item(219,26).Value2 = 0
$ws.Cells.Item(219,27).Value2 = 0
# row 220
$ws.Cells.Item(220,1).Value2 = 218
$ws.Cells.Item(220,2).Value2 = 7641720
$ws.Cells.Item(220,3).Value2 = "Mexico Liga de Expansion"
$ws.Cells.Item(220,4).Value2 = "Mexico Liga de Expansion"
$ws.Cells.Item(220,5).Value2 = 45386.00347222222
$ws.Cells.Item(220,6).Value2 = "Venados FC"
$ws.Cells.Item(220,7).Value2 = "Cancun FC"
$ws.Cells.Item(220,11).Value2 = 2.5
$ws.Cells.Item(220,12).Value2 = 3.1
$ws.Cells.Item(220,13).Value2 = 2.75
$ws.Cells.Item(220,14).Value2 = 2.55
$ws.Cells.Item(220,15).Value2 = 3.1
$ws.Cells.Item(220,16).Value2 = 2.7
$ws.Cells.Item(220,17).Value2 = 0
$ws.Cells.Item(220,18).Value2 = 1.85
$ws.Cells.Item(220,19).Value2 = 1.95
$ws.Cells.Item(220,20).Value2 = 2.5
$ws.Cells.Item(220,21).Value2 = 2
$ws.Cells.Item(220,22).Value2 = 1.8
$ws.Cells.Item(220,23).Value2 = 0
$ws.Cells.Item(220,24).Value2 = 0
$ws.Cells.Item(220,25).Value2 = 0
$ws.Cells.Item(220,26).Value2 = 0
$ws.Cells.Item(220,27).Value2 = 0
# row 221
$ws.Cells.Item(221,1).Value2 = 219
$ws.Cells.Item(221,2).Value2 = 7641722
$ws.Cells.Item(221,3).Value2 = "Mexico Liga de Expansion"
$ws.Cells.Item(221,4).Value2 = "Mexico Liga de Expansion"
$ws.Cells.Item(221,5).Value2 = 45387.00347222222
$ws.Cells.Item(221,6).Value2 = "Tepatitlan FC"
$ws.Cells.Item(221,7).Value2 = "Tlaxcala FC"
$ws.Cells.Item(221,11).Value2 = 2.6
$ws.Cells.Item(221,12).Value2 = 3.1
$ws.Cells.Item(221,13).Value2 = 2.65
$ws.Cells.Item(221,14).Value2 = 2.5
$ws.Cells.Item(221,15).Value2 = 3.1
$ws.Cells.Item(221,16).Value2 = 2.75
$ws.Cells.Item(221,17).Value2 = 0
$ws.Cells.Item(221,18).Value2 = 1.8
$ws.Cells.Item(221,19).Value2 = 2
$ws.Cells.Item(221,20).Value2 = 2.25
$ws.Cells.Item(221,21).Value2 = 1.85
$ws.Cells.Item(221,22).Value2 = 1.95
$ws.Cells.Item(221,23).Value2 = 0
$ws.Cells.Item(221,24).Value2 = 0
$ws.Cells.Item(221,25).Value2 = 0
$ws.Cells.Item(221,26).Value2 = 0
$ws.Cells.Item(221,27).Value2 = 0
